# timesheet-template.docx: turn the single-job {#jobs}...{/jobs} loop in
# the weekly-report table into a {#workBlocks}...{/workBlocks} loop that
# wraps just {jobId} ... {endTime}, and resize the table's grid columns
# to the new layout.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# --- resize the grid columns (w:tblGrid/w:gridCol, in twips -> points) ---
$newColWidthsTwips = @(3472, 1289, 1918, 1597, 1560, 1408, 1390, 1793)
for ($i = 0; $i -lt $newColWidthsTwips.Count; $i++) {
    $t.Columns.Item($i + 1).Width = $newColWidthsTwips[$i] / 20
}

# --- row 3 (the placeholder/template row) text edits ---

# Col 1 (Job Name/Address): was empty -> opens the workBlocks loop.
$c1 = $t.Cell(3, 1)
$c1.Range.Text = "{#workBlocks}"
$c1.Range.Font.Size = 14
$c1.Range.Font.SizeBi = 14

# Col 2 (Job #): drop the old per-job loop opener, keep just the field.
$d.Content.Find.Execute("{#jobs}{jobId}", $false, $false, $false, $false, `
    $false, $true, 1, $false, "{jobId}", 2)

# Col 5 (End Time): drop the old per-job loop closer entirely.
$d.Content.Find.Execute("{endTime}{/jobs}", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 2)

# Col 7 (was empty): now carries the End Time field.
$c7 = $t.Cell(3, 7)
$c7.Range.Text = "{endTime}"
$c7.Range.Font.Size = 14
$c7.Range.Font.SizeBi = 14

# Col 8 (Total Hours): was empty -> closes the workBlocks loop.
$c8 = $t.Cell(3, 8)
$c8.Range.Text = "{/workBlocks}"
$c8.Range.Font.Size = 14
$c8.Range.Font.SizeBi = 14
